$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.332.30'
$ws.Range("E2").Value = '  -0.96%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.584.25'
$ws.Range("E3").Value = '  -1.86%  '

$ws.Range("E4").Value = '  +0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '576.99'
$ws.Range("E5").Value = '  -2.73%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '192.43'
$ws.Range("E6").Value = '  +0.45%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.578.04'
$ws.Range("E7").Value = '  -1.83%  '

$ws.Range("E8").Value = '  -0.31%  '

$ws.Range("E9").Value = '  +0.01%  '

$ws.Range("E10").Value = '  -2.94%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.150'
$ws.Range("E11").Value = '  -0.52%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.50'
$ws.Range("E12").Value = '  -3.37%  '

$ws.Range("E13").Value = '  +1.40%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.93'
$ws.Range("E14").Value = '  -1.90%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.156.22'
$ws.Range("E15").Value = '  -1.90%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.582.52'
$ws.Range("E16").Value = '  -1.92%  '

$ws.Range("E17").Value = '  -0.81%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '67.334.54'
$ws.Range("E18").Value = '  -0.65%  '

$ws.Range("E19").Value = '  -0.70%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.35'
$ws.Range("E20").Value = '  -2.22%  '

$ws.Range("E21").Value = '  -2.64%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '401.61'
$ws.Range("E22").Value = '  +0.89%  '

$ws.Range("E23").Value = '  +21.86%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.19'
$ws.Range("E24").Value = '  -3.71%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.45'
$ws.Range("E25").Value = '  -2.21%  '

$ws.Range("E26").Value = '  +0.08%  '

$ws.Range("E27").Value = '  +0.80%  '

$ws.Range("E28").Value = '  +0.55%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.77'
$ws.Range("E29").Value = '  +3.58%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.07'
$ws.Range("E30").Value = '  +11.61%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '9.12'
$ws.Range("E31").Value = '  -1.67%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '31.24'
$ws.Range("E32").Value = '  -1.25%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '668.48'
$ws.Range("E33").Value = '  +10.69%  '

$ws.Range("E34").Value = '  -0.27%  '

$ws.Range("E35").Value = '  +1.12%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '64.00'
$ws.Range("E36").Value = '  -2.17%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '42.33'
$ws.Range("E37").Value = '  -3.75%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.420'
$ws.Range("E38").Value = '  +7.45%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.999'
$ws.Range("E39").Value = '  -0.13%  '

$ws.Range("E40").Value = '  +3.50%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.91'
$ws.Range("E41").Value = '  +15.76%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.12'
$ws.Range("E42").Value = '  +9.01%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.165.18'
$ws.Range("E43").Value = '  +14.56%  '

$ws.Range("E44").Value = '  -0.86%  '

$ws.Range("E45").Value = '  -0.25%  '

$ws.Range("E46").Value = '  -1.35%  '

$ws.Range("E47").Value = '  -2.37%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.12'
$ws.Range("E48").Value = '  +0.54%  '

$ws.Range("B49").Value = 'Monero'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '142.40'
$ws.Range("E49").Value = '  -0.86%  '

$ws.Range("B50").Value = 'THORChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.71'
$ws.Range("E50").Value = '  +0.48%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.51'
$ws.Range("E51").Value = '  -3.10%  '
